$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet logs availability checks in batches of 14 rows (one row per
# monitored service). Each run pushes the previous timestamps one batch
# down and stamps the newest batch (rows 2-15) with the current run time.

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44267.50927886293
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44267.48789518519
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44267.4664978125
}
